$d = $word.ActiveDocument

# --- Locate the paragraph containing "\begin{figure}[!htbp]" ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  $p = $d.Paragraphs.Item($i)
  if ($p.Range.Text.Contains("\begin{figure}[!htbp]")) {
    $target = $p
    break
  }
}

if ($target -eq $null) {
  Write-Output "ERROR: anchor paragraph not found"
} else {
  # Insert a brand new (empty) paragraph right after it. Using
  # InsertParagraphAfter keeps any bookmark that immediately follows the
  # anchor paragraph correctly positioned before the *next* paragraph
  # (here, the "cloud.sagemath" heading) rather than before our new one.
  $target.Range.InsertParagraphAfter()

  # Re-fetch the paragraph list; the newly inserted (still empty) paragraph
  # is immediately after the anchor paragraph.
  $newPara = $target.Next()
  $newStart = $newPara.Range.Start
  $newEnd = $newPara.Range.End

  $part1 = "Take a look at this writeLaTeX"
  $sp1 = " "
  $linkText = "http://goo.gl/k83ZHi"
  $sp2 = " "
  $part2 = "template to play around with this."

  $fullText = $part1 + $sp1 + $linkText + $sp2 + $part2

  # Replace the whole (still-empty) paragraph range -- including its
  # paragraph mark -- with clean OOXML for five plain runs. Rebuilding the
  # paragraph this way (rather than typing text into it) avoids inheriting
  # the "SourceCode"/"VerbatimChar" styling that the preceding paragraph
  # carries.
  $wholeParaRange = $d.Range($newStart, $newEnd)
  $xmlFragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p>' + `
    '<w:r><w:t xml:space="preserve">' + $part1 + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">' + $sp1 + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">' + $linkText + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">' + $sp2 + '</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve">' + $part2 + '</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
  $wholeParaRange.InsertXML($xmlFragment)

  # Work out where the URL text ended up so it can be turned into a real
  # hyperlink (matching the style used throughout the rest of the document).
  $linkStart = $newStart + $part1.Length + $sp1.Length
  $linkEnd = $linkStart + $linkText.Length
  $linkRange = $d.Range($linkStart, $linkEnd)

  $d.Hyperlinks.Add($linkRange, $linkText, [type]::Missing, [type]::Missing, [type]::Missing) | Out-Null

  Write-Output "Inserted writeLaTeX template paragraph with hyperlink."
}
